$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header columns B1/C1: DNNTSP <-> SFCNTSP
$ws.Range("B1").Value = "SFCNTSP"
$ws.Range("C1").Value = "DNNTSP"
# D1 stays "GP_top_freq" (unchanged)

# Row 2: weighted f1 (mean)
$ws.Range("A2").Value = "weighted f1(mean)"
$ws.Range("B2").Value = 0.4750733530506966
$ws.Range("C2").Value = 0.5059978957473863
$ws.Range("D2").Value = 0.5556464286094351

# Row 3: weighted f1 (std) -- new semantics for what used to be "weighted ROC AUC"
$ws.Range("A3").Value = "weighted f1(std)"
$ws.Range("B3").Value = 0.06345638024699643
$ws.Range("C3").Value = 0.008413442010530424
$ws.Range("D3").Value = 0

# Row 4: hamming_loss (mean)
$ws.Range("A4").Value = "hamming_loss(mean)"
$ws.Range("B4").Value = 0.106969696969697
$ws.Range("C4").Value = 0.1067454545454545
$ws.Range("D4").Value = 0.09909090909090909

# Row 5 (new): hamming_loss (std)
$ws.Range("A5").Value = "hamming_loss(std)"
$ws.Range("B5").Value = 0.003385541558408321
$ws.Range("C5").Value = 0.002326562391594835
$ws.Range("D5").Value = 0

# Row 6 (new): weighted ROC AUC (mean)
$ws.Range("A6").Value = "weighted ROC AUC(mean)"
$ws.Range("B6").Value = 0.6080632227869875
$ws.Range("C6").Value = 0.655800124848401
$ws.Range("D6").Value = 0.7092846855246081

# Row 7 (new): weighted ROC AUC (std)
$ws.Range("A7").Value = "weighted ROC AUC(std)"
$ws.Range("B7").Value = 0.06401437436429813
$ws.Range("C7").Value = 0.008144567490292331
$ws.Range("D7").Value = 0

# New A5:A7 cells should carry the same style as the other row-label cells (A2:A4)
$ws.Range("A2").Copy()
$ws.Range("A5:A7").PasteSpecial(-4122)
